$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 9: Distill, My Heart | Distilled Water
$ws_ALC.Range("H9").Value = 212.42857
$ws_ALC.Range("I9").Value = 245.36363
$ws_ALC.Range("K9").Value = 245.36363
$ws_ALC.Range("M9").Value = -76.36363

# ALC row 13: The Hexster Runoff | Maple Picatrix
$ws_ALC.Range("N13").ClearContents()
$ws_ALC.Range("H13").Value = 5
$ws_ALC.Range("I13").Value = 5
$ws_ALC.Range("J13").Value = 0
$ws_ALC.Range("K13").Value = 5
$ws_ALC.Range("L13").Value = 0
$ws_ALC.Range("M13").Value = 164

# ALC row 33: Glazed and Confused | Clear Glass Lens
$ws_ALC.Range("H33").Value = 280.66666
$ws_ALC.Range("I33").Value = 280.66666
$ws_ALC.Range("K33").Value = 280.66666
$ws_ALC.Range("M33").Value = -51.66665999999998

# ALC row 48: The Sting of Conscience | Sleeping Potion
$ws_ALC.Range("H48").Value = 1328
$ws_ALC.Range("J48").Value = 2649
$ws_ALC.Range("L48").Value = 7947
$ws_ALC.Range("N48").Value = -8531

# ALC row 56: Sleepless in Silvertear | Potent Sleeping Potion
$ws_ALC.Range("H56").Value = 1328
$ws_ALC.Range("J56").Value = 2649
$ws_ALC.Range("L56").Value = 7947
$ws_ALC.Range("N56").Value = -9015

# ALC row 115: 5-bell Energy | Competent Craftsman's Syrup
$ws_ALC.Range("H115").Value = 320
$ws_ALC.Range("I115").Value = 320
$ws_ALC.Range("K115").Value = 960
$ws_ALC.Range("M115").Value = 607

# ALC row 131: Mindful Study | Grade 5 Tincture of Mind
$ws_ALC.Range("H131").Value = 2297.1667
$ws_ALC.Range("I131").Value = 723.75
$ws_ALC.Range("J131").Value = 5444
$ws_ALC.Range("K131").Value = 2171.25
$ws_ALC.Range("L131").Value = 16332
$ws_ALC.Range("M131").Value = 2868.75
$ws_ALC.Range("N131").Value = -26412

# ALC row 132: Fast-forwarding Flora | Growth Formula Lambda
$ws_ALC.Range("H132").Value = 6668.3335
$ws_ALC.Range("I132").Value = 4334.6665
$ws_ALC.Range("K132").Value = 13003.9995
$ws_ALC.Range("M132").Value = -10473.9995

# ALC row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws_ALC.Range("H138").Value = 373406.34
$ws_ALC.Range("I138").Value = 1248.625
$ws_ALC.Range("J138").Value = 530104.3
$ws_ALC.Range("K138").Value = 3745.875
$ws_ALC.Range("L138").Value = 1590312.9
$ws_ALC.Range("M138").Value = 1394.125
$ws_ALC.Range("N138").Value = -1600592.9

# ARM row 6: Don't Hit Me One More Time | Bronze Hoplon
$ws_ARM.Range("H6").Value = 4007179.5
$ws_ARM.Range("I6").Value = 6671966.5
$ws_ARM.Range("K6").Value = 6671966.5
$ws_ARM.Range("M6").Value = -6671793.5

# ARM row 43: They've Got Legs | Steel Sabatons
$ws_ARM.Range("N43").ClearContents()
$ws_ARM.Range("H43").Value = 46342
$ws_ARM.Range("J43").Value = 0
$ws_ARM.Range("L43").Value = 0

# ARM row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws_ARM.Range("H61").Value = 3019.6
$ws_ARM.Range("I61").Value = 6000
$ws_ARM.Range("K61").Value = 6000
$ws_ARM.Range("M61").Value = -5788

# ARM row 63: Rivets Run through It | Mythrite Rivets
$ws_ARM.Range("N63").ClearContents()
$ws_ARM.Range("H63").Value = 2164.3333
$ws_ARM.Range("I63").Value = 2164.3333
$ws_ARM.Range("J63").Value = 0
$ws_ARM.Range("K63").Value = 2164.3333
$ws_ARM.Range("L63").Value = 0
$ws_ARM.Range("M63").Value = -1478.3333

# ARM row 66: A Riveting Revival (L) | Mythrite Rivets
$ws_ARM.Range("N66").ClearContents()
$ws_ARM.Range("H66").Value = 2164.3333
$ws_ARM.Range("I66").Value = 2164.3333
$ws_ARM.Range("J66").Value = 0
$ws_ARM.Range("K66").Value = 10821.6665
$ws_ARM.Range("L66").Value = 0
$ws_ARM.Range("M66").Value = -7389.666499999999

# ARM row 74: As the Bolt Flies | Titanium Nugget
$ws_ARM.Range("H74").Value = 1807.421
$ws_ARM.Range("I74").Value = 1463.0588
$ws_ARM.Range("J74").Value = 4734.5
$ws_ARM.Range("K74").Value = 1463.0588
$ws_ARM.Range("L74").Value = 4734.5
$ws_ARM.Range("M74").Value = -589.0588
$ws_ARM.Range("N74").Value = -6482.5

# ARM row 77: Heavy Metal Banned (L) | Titanium Nugget
$ws_ARM.Range("H77").Value = 1807.421
$ws_ARM.Range("I77").Value = 1463.0588
$ws_ARM.Range("J77").Value = 4734.5
$ws_ARM.Range("K77").Value = 7315.294
$ws_ARM.Range("L77").Value = 23672.5
$ws_ARM.Range("M77").Value = -2947.294
$ws_ARM.Range("N77").Value = -32408.5

# ARM row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws_ARM.Range("H132").Value = 2014.0667
$ws_ARM.Range("I132").Value = 2047
$ws_ARM.Range("J132").Value = 1800
$ws_ARM.Range("K132").Value = 6141
$ws_ARM.Range("L132").Value = 5400
$ws_ARM.Range("M132").Value = -3611
$ws_ARM.Range("N132").Value = -10460

# ARM row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws_ARM.Range("H136").Value = 3019.6
$ws_ARM.Range("I136").Value = 6000
$ws_ARM.Range("K136").Value = 18000
$ws_ARM.Range("M136").Value = -15450

# BSM row 5: Axe Me Anything | Bronze War Axe
$ws_BSM.Range("H5").Value = 195.33333
$ws_BSM.Range("I5").Value = 195.33333
$ws_BSM.Range("K5").Value = 195.33333
$ws_BSM.Range("M5").Value = -82.33332999999999

# BSM row 7: Thank You for Your Business | Bronze Bastard Sword
$ws_BSM.Range("N7").ClearContents()
$ws_BSM.Range("H7").Value = 10000150
$ws_BSM.Range("I7").Value = 10000150
$ws_BSM.Range("J7").Value = 0
$ws_BSM.Range("K7").Value = 10000150
$ws_BSM.Range("L7").Value = 0
$ws_BSM.Range("M7").Value = -10000037

# BSM row 102: Renting Mortality | Doman Steel Mortar
$ws_BSM.Range("H102").Value = 10496.5
$ws_BSM.Range("I102").Value = 10496.5
$ws_BSM.Range("K102").Value = 10496.5
$ws_BSM.Range("M102").Value = -7251.5

# CRP row 2: In with the New | Bone Harpoon
$ws_CRP.Range("H2").Value = 397.5
$ws_CRP.Range("J2").Value = 495
$ws_CRP.Range("L2").Value = 495
$ws_CRP.Range("N2").Value = -721

# CRP row 3: Touch and Heal | Maple Pattens
$ws_CRP.Range("H3").Value = 2828.5
$ws_CRP.Range("I3").Value = 2828.5
$ws_CRP.Range("K3").Value = 2828.5
$ws_CRP.Range("M3").Value = -2715.5

# CRP row 7: Gridania's Got Talent | Maple Lumber
$ws_CRP.Range("H7").Value = 350
$ws_CRP.Range("I7").Value = 350
$ws_CRP.Range("K7").Value = 350
$ws_CRP.Range("M7").Value = -237

# CRP row 13: Compulsory Conjury | Maple Cane
$ws_CRP.Range("H13").Value = 375
$ws_CRP.Range("J13").Value = 500
$ws_CRP.Range("L13").Value = 500
$ws_CRP.Range("N13").Value = -778

# CRP row 132: Hull Lotta Damage | Ginseng Lumber
$ws_CRP.Range("H132").Value = 1621.3636
$ws_CRP.Range("I132").Value = 1014.1667
$ws_CRP.Range("K132").Value = 3042.5001
$ws_CRP.Range("M132").Value = -512.5001000000002

# CUL row 5: What a Sap | Maple Syrup
$ws_CUL.Range("H5").Value = 2979.2
$ws_CUL.Range("J5").Value = 2999.75
$ws_CUL.Range("L5").Value = 8999.25
$ws_CUL.Range("N5").Value = -9223.25

# CUL row 63: The Next to Last Supper | Stuffed Cabbage Rolls
$ws_CUL.Range("M63").ClearContents()
$ws_CUL.Range("H63").Value = 0
$ws_CUL.Range("I63").Value = 0
$ws_CUL.Range("K63").Value = 0

# CUL row 66: Nostalgia through the Stomach (L) | Stuffed Cabbage Rolls
$ws_CUL.Range("M66").ClearContents()
$ws_CUL.Range("H66").Value = 0
$ws_CUL.Range("I66").Value = 0
$ws_CUL.Range("K66").Value = 0

# CUL row 117: A Good Omen | Peppered Popotoes
$ws_CUL.Range("H117").Value = 1300
$ws_CUL.Range("J117").Value = 1950
$ws_CUL.Range("L117").Value = 5850
$ws_CUL.Range("N117").Value = -12734

# CUL row 132: More Mezcal | Cooking Mezcal
$ws_CUL.Range("H132").Value = 1666.6666
$ws_CUL.Range("I132").Value = 1000
$ws_CUL.Range("J132").Value = 3000
$ws_CUL.Range("K132").Value = 9000
$ws_CUL.Range("L132").Value = 27000
$ws_CUL.Range("M132").Value = -6470
$ws_CUL.Range("N132").Value = -32060

# CUL row 135: Not-so-secret Ingredient | Royal Maple Syrup
$ws_CUL.Range("H135").Value = 2979.2
$ws_CUL.Range("J135").Value = 2999.75
$ws_CUL.Range("L135").Value = 26997.75
$ws_CUL.Range("N135").Value = -32067.75

# GSM row 20: Brothers in Arms | Brass Wristlets of Crafting
$ws_GSM.Range("H20").Value = 200
$ws_GSM.Range("I20").Value = 200
$ws_GSM.Range("K20").Value = 200
$ws_GSM.Range("M20").Value = 45

# GSM row 132: On Board for Lar | Lar Ingot
$ws_GSM.Range("H132").Value = 2918.7058
$ws_GSM.Range("I132").Value = 2508
$ws_GSM.Range("J132").Value = 5999
$ws_GSM.Range("K132").Value = 7524
$ws_GSM.Range("L132").Value = 17997
$ws_GSM.Range("M132").Value = -4994
$ws_GSM.Range("N132").Value = -23057

# LTW row 7: Tan Before the Ban | Leather
$ws_LTW.Range("H7").Value = 3783.6667
$ws_LTW.Range("I7").Value = 3787.2
$ws_LTW.Range("K7").Value = 3787.2
$ws_LTW.Range("M7").Value = -3675.2

# LTW row 46: Supply Side Logic | Boar Leather
$ws_LTW.Range("H46").Value = 821.75
$ws_LTW.Range("I46").Value = 762.3333
$ws_LTW.Range("K46").Value = 762.3333
$ws_LTW.Range("M46").Value = -574.3333

# LTW row 126: Battered Books | Saiga Leather
$ws_LTW.Range("H126").Value = 3783.6667
$ws_LTW.Range("I126").Value = 3787.2
$ws_LTW.Range("K126").Value = 11361.6
$ws_LTW.Range("M126").Value = -8891.599999999999

# LTW row 136: Respect for Br'aax | Br'aax Leather
$ws_LTW.Range("H136").Value = 5156.7144
$ws_LTW.Range("I136").Value = 5682.8335
$ws_LTW.Range("J136").Value = 2000
$ws_LTW.Range("K136").Value = 17048.5005
$ws_LTW.Range("L136").Value = 6000
$ws_LTW.Range("M136").Value = -14498.5005
$ws_LTW.Range("N136").Value = -11100

# WVR row 3: Trew Enough | Hempen Chausses
$ws_WVR.Range("H3").Value = 16001.5
$ws_WVR.Range("I3").Value = 12003
$ws_WVR.Range("K3").Value = 12003
$ws_WVR.Range("M3").Value = -11889

# WVR row 4: Not Cool Enough | Hempen Undershirt
$ws_WVR.Range("H4").Value = 3922.4443
$ws_WVR.Range("I4").Value = 852.4
$ws_WVR.Range("J4").Value = 7760
$ws_WVR.Range("K4").Value = 852.4
$ws_WVR.Range("L4").Value = 7760
$ws_WVR.Range("M4").Value = -739.4
$ws_WVR.Range("N4").Value = -7986

# WVR row 20: Read the Fine Print | Cotton Shepherd's Tunic
$ws_WVR.Range("N20").ClearContents()
$ws_WVR.Range("H20").Value = 1904.5
$ws_WVR.Range("I20").Value = 1904.5
$ws_WVR.Range("J20").Value = 0
$ws_WVR.Range("K20").Value = 1904.5
$ws_WVR.Range("L20").Value = 0
$ws_WVR.Range("M20").Value = -1664.5

# WVR row 22: Better Shroud than Sorry | Cotton Kurta
$ws_WVR.Range("H22").Value = 1076
$ws_WVR.Range("I22").Value = 614
$ws_WVR.Range("J22").Value = 2000
$ws_WVR.Range("K22").Value = 614
$ws_WVR.Range("L22").Value = 2000
$ws_WVR.Range("M22").Value = -321
$ws_WVR.Range("N22").Value = -2586

# WVR row 51: After the Smock-down | Linen Smock
$ws_WVR.Range("H51").Value = 13535
$ws_WVR.Range("I51").Value = 13535
$ws_WVR.Range("K51").Value = 13535
$ws_WVR.Range("M51").Value = -13025

# WVR row 52: Party Animals | Linen Deerstalker
$ws_WVR.Range("H52").Value = 4023808.5
$ws_WVR.Range("I52").Value = 5017510.5
$ws_WVR.Range("J52").Value = 49000
$ws_WVR.Range("K52").Value = 5017510.5
$ws_WVR.Range("L52").Value = 49000
$ws_WVR.Range("M52").Value = -5017284.5
$ws_WVR.Range("N52").Value = -49452

# WVR row 126: A Polished Purchase | Snow Linen
$ws_WVR.Range("H126").Value = 3109.1667
$ws_WVR.Range("I126").Value = 3224.5
$ws_WVR.Range("J126").Value = 3051.5
$ws_WVR.Range("K126").Value = 9673.5
$ws_WVR.Range("L126").Value = 9154.5
$ws_WVR.Range("M126").Value = -7203.5
$ws_WVR.Range("N126").Value = -14094.5

# WVR row 132: Comfy Cabins | Snow Cotton Cloth
$ws_WVR.Range("H132").Value = 3300.72
$ws_WVR.Range("I132").Value = 2069.1365
$ws_WVR.Range("K132").Value = 6207.4095
$ws_WVR.Range("M132").Value = -3677.4095
